# DEV-5077: apply new permissions terminology to the excel2xml test workbook
# "open" -> "public", "restricted" -> "private" (whole-cell matches only)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlWhole = 1 -> only replace cells whose entire content is "open"/"restricted"
$ws.Cells.Replace("open", "public", 1)
$ws.Cells.Replace("restricted", "private", 1)

# Reflect the editor's cursor/scroll position at the time of the edit:
# scrolled right so column D is left-most visible, with O18 as the active cell.
$ws.Range("O18").Select()
$excel.ActiveWindow.ScrollColumn = 4
